$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("historico")

# New answered-question rows to append after the last existing row (169).
$rows = @(
    @(183, '01/01/2026 04:17', 1284, 'Matemática Financeira', 'Juros Compostos', 0),
    @(184, '01/01/2026 04:28', 1257, 'Matemática Financeira', 'Juros Simples', 1),
    @(185, '01/01/2026 04:35', 1267, 'Matemática Financeira', 'Juros Compostos', 1),
    @(186, '01/01/2026 05:08', 1283, 'Matemática Financeira', 'Juros Compostos', 0),
    @(187, '01/01/2026 16:10', 1715, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(188, '01/01/2026 16:11', 1721, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(189, '01/01/2026 16:11', 1714, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(190, '01/01/2026 16:12', 1719, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(191, '01/01/2026 16:14', 1718, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(192, '01/01/2026 16:18', 1723, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(193, '01/01/2026 16:22', 1724, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(194, '01/01/2026 16:22', 1722, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(195, '01/01/2026 16:23', 1716, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 1),
    @(196, '01/01/2026 16:24', 1717, 'Contabilidade De Custos', 'Custos: Conceitos Introdutórios', 0)
)

$startRow = 170
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
